# Updated the workspace name variables.
#
# Sheet1 ("_Test_Suite_Statistics") rows 2-5 hold one "workspace name"
# record per row (name / total cases / automated cases / status / note).
# The rows get re-sorted alphabetically by column A, which cyclically
# shifts the data (and the trailing note in column E) up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: EnableScheduler -----------------------------------------
$ws.Range("A2").Value = "EnableScheduler"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = "Automated"
$ws.Range("E2").Clear()

# --- Row 3: ProgressBar ----------------------------------------------
$ws.Range("A3").Value = "ProgressBar"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "Suited to Manual"

# --- Row 4: RequiredFields --------------------------------------------
$ws.Range("A4").Value = "RequiredFields"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "Automated"

# --- Row 5: Step1Fields ------------------------------------------------
$ws.Range("A5").Value = "Step1Fields"
$ws.Range("B5").Value = 19
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = "Automated"
$ws.Range("E5").Value = "One test case is untested and waiting on TBH143"

# Re-apply the sheet's sort (column A, A2:E5) so the sortState the
# range was already sorted under is refreshed to the shrunk range.
$sortRange = $ws.Range("A2:E5")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2")) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Move the selection cursor like the saved workbook shows.
$ws.Range("A10").Select() | Out-Null
